$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = 'Each'

$ws.Range("C8").Value = 7

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.0'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'

$ws.Range("F8").Value = 50

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '350.00'

$ws.Range("C9").Value = 2

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.0'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'

$ws.Range("F9").Value = 30

$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '60.00'

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = 'Each'

$ws.Range("C10").Value = 75

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '8.0'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'

$ws.Range("F10").Value = 30

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '2250.00'

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = 'Each'

$ws.Range("C11").Value = 52

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '10.0'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'

$ws.Range("F11").Value = 303

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '15756.00'

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = 'R. mtr.'

$ws.Range("C12").Value = 9

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '20 mm'

$ws.Range("F12").Value = 40

$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '360.00'

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = 'Mtr.'

$ws.Range("C13").Value = 12

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'

$ws.Range("F13").Value = 122

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '1464.00'

$ws.Range("C14").Value = 6

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'

$ws.Range("A15").Value = ""

$ws.Range("C15").Value = 62

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '36'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = 'Total'

$ws.Range("F15").Value = 0

$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '0.00'

$ws.Range("A16").Value = ""

$ws.Range("C16").Value = 100

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '38'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = 'Grand Total'

$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '20240.00'

$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = '20240.00'

$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '20240.00'

$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '20240.00'
